$wb = $excel.ActiveWorkbook

# --- Parameters sheet: replace MCMC-related parameters with the new
#     max_spectra_per_peptide / nbatch pair, and drop the now-unused rows.
$ws = $wb.Worksheets.Item("Parameters")

$ws.Range("A2").Value = "max_spectra_per_peptide"
$ws.Range("B2").Value = 5
$ws.Range("A3").Value = "nbatch"
$ws.Range("B3").Value = 100

# Remove the old rows 4-11 (model_fc through nworker no longer exist).
$ws.Rows("4:11").Delete()

# Make Parameters the active/selected sheet, with F10 as the active cell.
$ws.Activate()
$ws.Range("F10").Select() | Out-Null
